# Append a new day's row (2025-04-12) to every price sheet in the
# workbook, carrying forward the same price that was recorded on the
# previous day (2025-04-11, row 41). Values are written as literal text
# (via a leading apostrophe) so they stay plain strings -- matching the
# existing "Date"/"Price" columns, which are stored as text rather than
# real Excel dates/numbers -- instead of being auto-converted to a date
# serial number or a number.

$wb = $excel.ActiveWorkbook

$sheetPrices = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "41.5"
    "N-type Wafer"               = "1.25"
    "Cell Topcon 183mm"          = "0.303"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,192"
    "Silver Busbar front-side"   = "7,773"
    "Silver finger front-side"   = "7,823"
    "USD_CNY"                    = "7.3258"
}

foreach ($sheetName in $sheetPrices.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $price = $sheetPrices[$sheetName]

    $ws.Cells.Item(42, 1).Value = "'2025-04-12"
    $ws.Cells.Item(42, 2).Value = "'" + $price
}
